# adding year from admin Type: SAVE.
# Shifts the FICA rate table on the "Configs" sheet down by one year
# (a new, most-recent year block is introduced and the oldest year's
# data, previously on rows 17-19, is duplicated onto new rows 20-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configs")

# --- Header cell B2: update the formula-name label text -------------
$ws.Range("B2").Value = "SmartRules PaymentToolsFicaDetails FicaTaxRates(PaymentToolsFicaInput input)"

# --- Shift each 3-row "year block" (rows 5-19) forward by one year --
# Only the top (anchor) row of every merged B5:B7 / B8:B10 / ... block
# can be written to directly; the other two rows are merged into it.
$ws.Range("B5").Value  = 2026
$ws.Range("B8").Value  = 2025
$ws.Range("B11").Value = 2024
$ws.Range("B14").Value = 2023
$ws.Range("B17").Value = 2022

# --- FSST threshold (column E) values cascade down from the row below
$ws.Range("E13").Value = 168600
$ws.Range("E16").Value = 152000
$ws.Range("E19").Value = 147000

# --- Row 19's "Rate" cell was stored as text "6.2"; make it numeric --
$ws.Range("D19").Value = 6.2

# --- Add the new oldest-year (2021) block on rows 20-22, matching the
#     layout/values that used to describe 2021 further up the table --
$ws.Range("B20").Value = 2021
$ws.Range("C20").Value = "AFMT"
$ws.Range("D20").Value = 0.9
$ws.Range("E20").Value = 200000

$ws.Range("B21").Value = 2021
$ws.Range("C21").Value = "FMT"
$ws.Range("D21").Value = 1.45

$ws.Range("B22").Value = 2021
$ws.Range("C22").Value = "FSST"
$ws.Range("D22").Value = 6.2
$ws.Range("E22").Value = 142800

# --- New trailing blank separator row ---------------------------------
# (write then clear so the cell/row exists and inherits the default
#  column style, matching the blank separator rows used elsewhere)
$ws.Range("B23").Value = 0
$ws.Range("B23").ClearContents()

Write-Host "Configs sheet year rollover applied."
